# Add new tracker abbreviations (ACC_X / ACC_Y / ACC_Z) to the datapoints
# section of the tracker import-specification sheet. These three new
# datapoint rows are inserted immediately above the existing AZIMUTH row,
# pushing AZIMUTH..QS_RX down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "datapoints" block currently starts at row 12 with AZIMUTH. Insert
# three blank rows above it so the new ACC_* entries can be written in,
# and everything below (AZIMUTH ... QS_RX) shifts down to rows 15-27.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Row 12: ACC_X
$ws.Range("A12").Value = "datapoints"
$ws.Range("B12").Value = "ACC_X"
$ws.Range("C12").Value = "m/s²"
$ws.Range("D12").Value = "Sensor acceleration on east-west axis"

# Row 13: ACC_Y
$ws.Range("A13").Value = "datapoints"
$ws.Range("B13").Value = "ACC_Y"
$ws.Range("C13").Value = "m/s²"
$ws.Range("D13").Value = "Sensor acceleration on north-south axis"

# Row 14: ACC_Z
$ws.Range("A14").Value = "datapoints"
$ws.Range("B14").Value = "ACC_Z"
$ws.Range("C14").Value = "m/s²"
$ws.Range("D14").Value = "Sensor acceleration on vertical axis"

Write-Output "Inserted ACC_X / ACC_Y / ACC_Z datapoint rows"
